# Generate Report for Handoff
#
# Updates the localization-status report:
#  - "Overview" sheet: refresh the "Latest HO Xliff Generate Date" timestamp
#    for the rows tied to the newly generated handoff (rows 7,8,9,11,13,14).
#  - "zh-cn" / "de-de" sheets: refresh the matching "Latest Handoff Datetime"
#    timestamp for those same rows, and mark their "Priority" column as "ht"
#    (handoff type) now that a handoff file has been produced for them.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# Overview sheet - column G is "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-02 12:25:01"
}

# zh-cn sheet - column H is "Latest Handoff Datetime", column E is "Priority"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-09-02 12:24:56"
    $wsZhCn.Range("E$r").Value = "ht"
}

# de-de sheet - column H is "Latest Handoff Datetime", column E is "Priority"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-09-02 12:25:01"
    $wsDeDe.Range("E$r").Value = "ht"
}
